# Sync automático del tracker (cada 3h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Delete the two trailing duplicate rows (7 and 8) ---
$ws.Rows("7:8").Delete()

# --- Re-key the event_id column (A) as text for the remaining data rows ---
# Excel stores a leading apostrophe as a "text-quoted number" (quotePrefix),
# keeping the digits intact while forcing a text cell type instead of numeric.
$ws.Cells.Item(2,1).Value = "'14494923"
$ws.Cells.Item(3,1).Value = "'14494919"
$ws.Cells.Item(4,1).Value = "'14494979"
$ws.Cells.Item(6,1).Value = "'14578002"

# --- Row 5 is replaced with a new match entirely ---
$ws.Cells.Item(5,1).Value = "'14487604"
$ws.Cells.Item(5,2).Value = "'2025-08-30"
$ws.Cells.Item(5,3).Value = "Daniel Rincon"
$ws.Cells.Item(5,4).Value = "Harold Mayot"
$ws.Cells.Item(5,5).Value = "Gana Daniel Rincon"
$ws.Cells.Item(5,6).Value = 2.75
